$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 field updates (per commit diff) ---

# DATE_TYPE_CODE: 001 -> 002 (keep as text, preserve leading zero)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"
$ws.Range("J2").Style = "Normal"

# NOTICE_DATE / REPORT_DATE
$ws.Range("M2").Value = "2020-12-24 00:00:00"
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Numeric financial figures
$ws.Range("O2").Value = 562078648.1
$ws.Range("P2").Value = 87549711.44
$ws.Range("Q2").Value = 64227611.37

# MONETARYFUNDS_RATIO cleared
$ws.Range("R2").ClearContents()

$ws.Range("S2").Value = 160109216.78

# ACCOUNTS_RECE_RATIO cleared
$ws.Range("T2").ClearContents()

$ws.Range("U2").Value = 128079354.25

# INVENTORY_RATIO cleared
$ws.Range("V2").ClearContents()

$ws.Range("W2").Value = 254762296.94
$ws.Range("X2").Value = 146302260.58

# ACCOUNTS_PAYABLE_RATIO cleared
$ws.Range("Y2").ClearContents()

$ws.Range("Z2").Value = 206026.61

# ADVANCE_RECEIVABLES_RATIO cleared
$ws.Range("AA2").ClearContents()

$ws.Range("AB2").Value = 307316351.16

# TOTAL_EQUITY_RATIO / TOTAL_ASSETS_RATIO / TOTAL_LIAB_RATIO cleared
$ws.Range("AC2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()

$ws.Range("AF2").Value = 188.8014555444
$ws.Range("AG2").Value = 45.3250266313
